$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 53641.438
$ws.Range("J17").Value = 53641.438
$ws.Range("L17").Value = 160924.314
$ws.Range("N17").Value = -161260.314

$ws.Range("H30").Value = 5500
$ws.Range("J30").Value = 5500
$ws.Range("L30").Value = 16500
$ws.Range("N30").Value = -16702

$ws.Range("H42").Value = 150
$ws.Range("I42").Value = 150
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 450
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -220
$ws.Range("N42").ClearContents()

$ws.Range("H53").Value = 415.03705
$ws.Range("I53").Value = 553.375
$ws.Range("J53").Value = 356.78946
$ws.Range("K53").Value = 553.375
$ws.Range("L53").Value = 356.78946
$ws.Range("M53").Value = 83.625
$ws.Range("N53").Value = -1630.78946

$ws.Range("H82").Value = 3488.4
$ws.Range("I82").Value = 480.66666
$ws.Range("J82").Value = 8000
$ws.Range("K82").Value = 1441.99998
$ws.Range("L82").Value = 24000
$ws.Range("M82").Value = -1035.99998
$ws.Range("N82").Value = -24812

$ws.Range("H85").Value = 3488.4
$ws.Range("I85").Value = 480.66666
$ws.Range("J85").Value = 8000
$ws.Range("K85").Value = 1441.99998
$ws.Range("L85").Value = 24000
$ws.Range("M85").Value = -37.99998000000005
$ws.Range("N85").Value = -26808

$ws.Range("H125").Value = 1684.7273
$ws.Range("J125").Value = 2024
$ws.Range("L125").Value = 18216
$ws.Range("N125").Value = -23136

$ws.Range("H132").Value = 4769056.5
$ws.Range("I132").Value = 6068808.5
$ws.Range("J132").Value = 3300
$ws.Range("K132").Value = 18206425.5
$ws.Range("L132").Value = 9900
$ws.Range("M132").Value = -18203895.5
$ws.Range("N132").Value = -14960

$ws.Range("H137").Value = 2737.2537
$ws.Range("I137").Value = 2766.6597
$ws.Range("J137").Value = 2668.15
$ws.Range("K137").Value = 8299.9791
$ws.Range("L137").Value = 8004.450000000001
$ws.Range("M137").Value = -5749.9791
$ws.Range("N137").Value = -13104.45

$ws.Range("H141").Value = 714714.25
$ws.Range("I141").Value = 2206.25
$ws.Range("J141").Value = 2139730.2
$ws.Range("K141").Value = 6618.75
$ws.Range("L141").Value = 6419190.600000001
$ws.Range("M141").Value = -1438.75
$ws.Range("N141").Value = -6429550.600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4953.6
$ws.Range("I32").Value = 4829.951
$ws.Range("K32").Value = 4829.951
$ws.Range("M32").Value = -4542.951

$ws.Range("H45").Value = 1618.4138
$ws.Range("I45").Value = 1093.1364
$ws.Range("K45").Value = 1093.1364
$ws.Range("M45").Value = -716.1364000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 300.2
$ws.Range("I22").Value = 300.2
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 300.2
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -127.2
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2652.25
$ws.Range("I132").Value = 2273.25
$ws.Range("J132").Value = 3599.75
$ws.Range("K132").Value = 6819.75
$ws.Range("L132").Value = 10799.25
$ws.Range("M132").Value = -4289.75
$ws.Range("N132").Value = -15859.25

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1204.7368
$ws.Range("I107").Value = 405.14285
$ws.Range("J107").Value = 1671.1666
$ws.Range("K107").Value = 1215.42855
$ws.Range("L107").Value = 5013.4998
$ws.Range("M107").Value = 704.5714499999999
$ws.Range("N107").Value = -8853.4998

$ws.Range("H132").Value = 1563
$ws.Range("I132").Value = 1206
$ws.Range("J132").Value = 2118.3333
$ws.Range("K132").Value = 10854
$ws.Range("L132").Value = 19064.9997
$ws.Range("M132").Value = -8324
$ws.Range("N132").Value = -24124.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3260.3333
$ws.Range("J126").Value = 3858.6667
$ws.Range("L126").Value = 11576.0001
$ws.Range("N126").Value = -16516.0001

$ws.Range("H132").Value = 2781.5
$ws.Range("I132").Value = 2541.697
$ws.Range("J132").Value = 3125.5652
$ws.Range("K132").Value = 7625.091
$ws.Range("L132").Value = 9376.695599999999
$ws.Range("M132").Value = -5095.091
$ws.Range("N132").Value = -14436.6956

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3848075.5
$ws.Range("I7").Value = 7693648.5
$ws.Range("J7").Value = 2502.1538
$ws.Range("K7").Value = 7693648.5
$ws.Range("L7").Value = 2502.1538
$ws.Range("M7").Value = -7693536.5
$ws.Range("N7").Value = -2726.1538

$ws.Range("H109").Value = 25700
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 25700
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 25700
$ws.Range("N109").Value = -28474
$ws.Range("M109").ClearContents()

$ws.Range("H126").Value = 3848075.5
$ws.Range("I126").Value = 7693648.5
$ws.Range("J126").Value = 2502.1538
$ws.Range("K126").Value = 23080945.5
$ws.Range("L126").Value = 7506.4614
$ws.Range("M126").Value = -23078475.5
$ws.Range("N126").Value = -12446.4614

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1171.5
$ws.Range("I81").Value = 1107.7693
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 2215.5386
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -1154.5386
$ws.Range("N81").Value = -6122

$ws.Range("H84").Value = 1171.5
$ws.Range("I84").Value = 1107.7693
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 11077.693
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -5773.692999999999
$ws.Range("N84").Value = -30608

$ws.Range("H126").Value = 2859232.5
$ws.Range("I126").Value = 1387.5333
$ws.Range("J126").Value = 20006302
$ws.Range("K126").Value = 4162.5999
$ws.Range("L126").Value = 60018906
$ws.Range("M126").Value = -1692.5999
$ws.Range("N126").Value = -60023846

$ws.Range("H132").Value = 31334
$ws.Range("I132").Value = 7044.75
$ws.Range("J132").Value = 48998.91
$ws.Range("K132").Value = 21134.25
$ws.Range("L132").Value = 146996.73
$ws.Range("M132").Value = -18604.25
$ws.Range("N132").Value = -152056.73
